$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(0.00608883429203496, -0.0175131295736708),
    @(-0.000576416619771311, 0.0727627847056778),
    @(-0.0031206859964281, 0.29511000506986),
    @(0.0070188703343616, -0.0460215205216533),
    @(-0.000930450724035223, 0.0914160872067665),
    @(-0.00371074193857551, 0.354864575417157),
    @(0.00773514152673103, -0.0915815999898393),
    @(-0.000583936216418073, 0.0783507720320212),
    @(-0.00342901587387804, 0.33142131727591),
    @(0.00664264299264749, -0.0351847653972643),
    @(-0.000692440472029338, 0.0729626400416282),
    @(-0.00285691633938786, 0.288228479195512),
    @(-0.000509580612836343, 0.0501315766734615),
    @(0.00614888577083178, -0.0492486139661487)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
